# Add 150% "line spacing within" (a:lnSpc/a:spcPct val="150000") to every
# paragraph of the bulleted/numbered list body placeholder on slide 3
# ("Types of chatbot").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The numbered list lives in the body placeholder, which is the 2nd shape
# on this slide (1st shape is the title).
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    # SpaceWithin is expressed as a multiplier of single spacing when the
    # line rule is percentage based (1.5 => 150%), which serializes to
    # <a:lnSpc><a:spcPct val="150000"/></a:lnSpc> in the paragraph's pPr.
    $para.ParagraphFormat.SpaceWithin = 1.5
}

Write-Host "Set 150% line spacing on $count paragraphs in slide 3's body placeholder"
